$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.644726333333334
$ws.Range("H2").Value = 4.934179
$ws.Range("I2").Value = 0.03084360558270512
$ws.Range("J2").Value = 0.03084360558270512
$ws.Range("M2").Value = 4.331589999999999
$ws.Range("N2").Value = 12.99477
$ws.Range("O2").Value = 0.1478799966101367
$ws.Range("P2").Value = 0.1478799966101367
$ws.Range("Q2").Value = 7.124280138203333
$ws.Range("R2").Value = 64.11852124383
$ws.Range("S2").Value = 0.004561152289014827
$ws.Range("T2").Value = 0.004561152289014827

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.644726333333334
$ws.Range("H3").Value = 4.934179
$ws.Range("I3").Value = 0.03084360558270512
$ws.Range("J3").Value = 0.03084360558270512
$ws.Range("O3").Value = 0.5404313285772905
$ws.Range("P3").Value = 0.5404313285772904
$ws.Range("Q3").Value = 26.03586873481245
$ws.Range("R3").Value = 234.322818613312
$ws.Range("S3").Value = 0.01666885074317526
$ws.Range("T3").Value = 0.01666885074317526

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.644726333333334
$ws.Range("H4").Value = 4.934179
$ws.Range("I4").Value = 0.03084360558270512
$ws.Range("J4").Value = 0.03084360558270512
$ws.Range("M4").Value = 9.129751000000001
$ws.Range("O4").Value = 0.3116886748125729
$ws.Range("P4").Value = 0.3116886748125729
$ws.Range("Q4").Value = 15.01594188647634
$ws.Range("R4").Value = 135.143476978287
$ws.Range("S4").Value = 0.009613602550515033
$ws.Range("T4").Value = 0.009613602550515033

# Row 5
$ws.Range("I5").Value = 0.828024694817689
$ws.Range("J5").Value = 0.828024694817689
$ws.Range("M5").Value = 4.331589999999999
$ws.Range("N5").Value = 12.99477
$ws.Range("O5").Value = 0.1478799966101367
$ws.Range("P5").Value = 0.1478799966101367
$ws.Range("Q5").Value = 191.25779155143
$ws.Range("R5").Value = 1721.32012396287
$ws.Range("S5").Value = 0.1224482890627494
$ws.Range("T5").Value = 0.1224482890627494

# Row 6
$ws.Range("I6").Value = 0.828024694817689
$ws.Range("J6").Value = 0.828024694817689
$ws.Range("O6").Value = 0.5404313285772905
$ws.Range("P6").Value = 0.5404313285772904
$ws.Range("S6").Value = 0.4474904859151292
$ws.Range("T6").Value = 0.4474904859151291

# Row 7
$ws.Range("I7").Value = 0.828024694817689
$ws.Range("J7").Value = 0.828024694817689
$ws.Range("M7").Value = 9.129751000000001
$ws.Range("O7").Value = 0.3116886748125729
$ws.Range("P7").Value = 0.3116886748125729
$ws.Range("Q7").Value = 403.1166416199271
$ws.Range("S7").Value = 0.2580859198398106
$ws.Range("T7").Value = 0.2580859198398106

# Row 8
$ws.Range("I8").Value = 0.1411316995996059
$ws.Range("J8").Value = 0.1411316995996059
$ws.Range("M8").Value = 4.331589999999999
$ws.Range("N8").Value = 12.99477
$ws.Range("O8").Value = 0.1478799966101367
$ws.Range("P8").Value = 0.1478799966101367
$ws.Range("Q8").Value = 32.59871034312999
$ws.Range("R8").Value = 293.38839308817
$ws.Range("S8").Value = 0.02087055525837255
$ws.Range("T8").Value = 0.02087055525837255

# Row 9
$ws.Range("I9").Value = 0.1411316995996059
$ws.Range("J9").Value = 0.1411316995996059
$ws.Range("O9").Value = 0.5404313285772905
$ws.Range("P9").Value = 0.5404313285772904
$ws.Range("S9").Value = 0.07627199191898607
$ws.Range("T9").Value = 0.07627199191898605

# Row 10
$ws.Range("I10").Value = 0.1411316995996059
$ws.Range("J10").Value = 0.1411316995996059
$ws.Range("M10").Value = 9.129751000000001
$ws.Range("O10").Value = 0.3116886748125729
$ws.Range("P10").Value = 0.3116886748125729
$ws.Range("Q10").Value = 68.70874398405701
$ws.Range("S10").Value = 0.04398915242224728
$ws.Range("T10").Value = 0.04398915242224728
